# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45190 (2023-09-21) to 45192 (2023-09-23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value = 45192
    }
}
